# Update TrialsSetup 2026-01-26 12:00
# Refresh of the "Days remaining" figures for two trials in the Query1
# table on Sheet1 (REJOICE and REMASTER), reflecting newer data pulled
# from the linked Power Query source.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 -> REJOICE (MK-5909-003): Days remaining 27 -> 24
$ws.Range("B11").Value = 24

# Row 14 -> REMASTER (CLOU): Days remaining 47 -> 44
$ws.Range("B14").Value = 44

$wb.Save()
